$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B2").Value = 0.3683274021352313
$ws.Range("C2").Value = 0.0731070496083551
$ws.Range("E2").Value = 0.1362530413625304
$ws.Range("F2").Value = 0.2828282828282828
$ws.Range("G2").Value = 0.6722068328716528
$ws.Range("H2").Value = 0.8006286784376673
$ws.Range("J2").Value = 355
$ws.Range("K2").Value = 179

# ---- Classification Report sheet ----
$ws = $wb.Worksheets.Item("Classification Report")
$ws.Range("C2").Value = 0.3352059925093633
$ws.Range("D2").Value = 0.5021037868162693

$ws.Range("B3").Value = 0.0731070496083551
$ws.Range("D3").Value = 0.1362530413625304

$ws.Range("B4").Value = 0.3683274021352313
$ws.Range("C4").Value = 0.3683274021352313
$ws.Range("D4").Value = 0.3683274021352313
$ws.Range("E4").Value = 0.3683274021352313

$ws.Range("B5").Value = 0.5365535248041775
$ws.Range("C5").Value = 0.6676029962546817
$ws.Range("D5").Value = 0.3191784140893998

$ws.Range("B6").Value = 0.9538202800516618
$ws.Range("C6").Value = 0.3683274021352313
$ws.Range("D6").Value = 0.4838763475409941

# ---- Confusion Matrix sheet ----
$ws = $wb.Worksheets.Item("Confusion Matrix")
$ws.Range("B2").Value = 179
$ws.Range("C2").Value = 355
